# Edit the "Rectangle 88" textbox on slide 2 of FlowChart.pptx:
#   - reposition/resize it (a:off / a:ext)
#   - update the formula text it contains, keeping the existing
#     run/formatting boundaries intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("Rectangle 88")

# --- Reposition / resize -------------------------------------------------
# Target EMU values (from the canonical OOXML):
#   off  x=8786956  y=9121107
#   ext  cx=4060727 cy=246221 (cy unchanged)
# PowerPoint's COM Left/Top/Width/Height are in points (1 pt = 12700 EMU)
# and are stored as single-precision floats, so the literals below are
# chosen so that they round-trip to the exact target EMU values.
$shp.Left = 691.8863525390625
$shp.Top = 718.1974487304688
$shp.Width = 319.7423095703125

# --- Update text content ---------------------------------------------------
$tr = $shp.TextFrame.TextRange

# Run 1: "ACTUAL_SCALE = DISTANCE_REF /" -> "ACTUAL_SCALE = DISTANCE_REF_MM /"
$newRun1Text = "ACTUAL_SCALE = DISTANCE_REF_MM /"
$run1 = $tr.Characters(1, 29)
$run1.Text = $newRun1Text

# Run 2 (single space separator) is untouched.
# Run 3: "835.7672/96 x 25.4" -> "835.7672 /96 /25.4"
$run3Start = $newRun1Text.Length + 1 + 1
$run3 = $shp.TextFrame.TextRange.Characters($run3Start, 18)
$run3.Text = "835.7672 /96 /25.4"
